$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187, shifting existing rows 187-262 down to 188-263
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new record's data
$ws.Range("A187").Value = 6
$ws.Range("B187").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C187").Value = 'Metropolitana'
$ws.Range("D187").Value = 44875
$ws.Range("E187").Value = 13
$ws.Range("F187").Value = 100112001
$ws.Range("G187").Value = 'Berenjena'
$ws.Range("H187").Value = 'Sin especificar'
$ws.Range("I187").Value = 'Primera'
$ws.Range("J187").Value = 400
$ws.Range("K187").Value = 17000
$ws.Range("L187").Value = 18000
$ws.Range("M187").Value = 17425
$ws.Range("N187").Value = '$/caja 50 unidades'
$ws.Range("O187").Value = 'Provincia de Huasco'
$ws.Range("P187").Value = 348
$ws.Range("Q187").Value = 50
$ws.Range("R187").Value = 'Hortaliza'
